# Update from Github Action
# Re-order the I / R / S columns (公积金 / 公积金比例 / 试用期工资 related data)
# and bump the "更新时间" timestamp on row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colI = 9   # I - 公积金 / 公积金比例
$colR = 18  # R - 公积金比例 / 试用期工资
$colS = 19  # S - 试用期工资 / 公积金

# Header row (1) and every data row except row 2 get their I/R/S values
# rotated: new I <- old R, new R <- old S, new S <- old I.
for ($r = 1; $r -le 20; $r++) {
    if ($r -eq 2) { continue }

    $i = $ws.Cells.Item($r, $colI).Value2
    $rVal = $ws.Cells.Item($r, $colR).Value2
    $s = $ws.Cells.Item($r, $colS).Value2

    $targets = @($colI, $colR, $colS)
    $newVals = @($rVal, $s, $i)

    for ($k = 0; $k -lt 3; $k++) {
        $cell = $ws.Cells.Item($r, $targets[$k])
        $val = $newVals[$k]

        # Plain-looking percentages (e.g. "100%", "8%") get auto-coerced into a
        # numeric percent value by Excel's normal text entry parsing. Force
        # those through as literal text (quote-prefix), then drop the
        # resulting quote-prefix style so the cell's formatting is untouched.
        if ($val -ne $null -and $val -is [string] -and $val -match '^\d+(\.\d+)?%$') {
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}

# Row 2's timestamp (column Q = 17) was refreshed by the sync action.
$ws.Cells.Item(2, 17).Value = "2022-02-10 07:19:05"
